# Auto-generated update: resum_diari_meteocat.xlsx
# Commit: Update automàtic: dades i banners [2026-02-07 04:19]
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-07 04:17:38'
$ws.Range('H2').NumberFormat = "@"
$ws.Range('H2').Value = '93%'
$ws.Range('N2').Value = '-1.6 °C 3:58 TU'
$ws.Range('O2').Value = '-1.2 °C'
$ws.Range('E3').Value = '2026-02-07 04:17:40'
$ws.Range('H3').NumberFormat = "@"
$ws.Range('H3').Value = '95%'
$ws.Range('O3').Value = '-5.5 °C'
$ws.Range('E4').Value = '2026-02-07 04:17:42'
$ws.Range('K4').Value = '-0.1 MJ/m2'
$ws.Range('L4').Value = '41.4 km/h - 313º 3:40 TU'
$ws.Range('O4').Value = '11.7 °C'
$ws.Range('E5').Value = '2026-02-07 04:17:45'
$ws.Range('J5').Value = '1000.9 hPa'
$ws.Range('E6').Value = '2026-02-07 04:17:47'
$ws.Range('H6').NumberFormat = "@"
$ws.Range('H6').Value = '59%'
$ws.Range('J6').Value = '1002.5 hPa'
$ws.Range('N6').Value = '11.3 °C 3:56 TU'
$ws.Range('O6').Value = '11.9 °C'
$ws.Range('E7').Value = '2026-02-07 04:17:50'
$ws.Range('H7').NumberFormat = "@"
$ws.Range('H7').Value = '75%'
$ws.Range('N7').Value = '7.0 °C 3:49 TU'
$ws.Range('O7').Value = '7.9 °C'
$ws.Range('E8').Value = '2026-02-07 04:17:52'
$ws.Range('N8').Value = '2.8 °C 3:59 TU'
$ws.Range('O8').Value = '4.5 °C'
$ws.Range('E9').Value = '2026-02-07 04:17:54'
$ws.Range('N9').Value = '0.5 °C 3:56 TU'
$ws.Range('O9').Value = '2.3 °C'
$ws.Range('E10').Value = '2026-02-07 04:17:57'
$ws.Range('E11').Value = '2026-02-07 04:17:59'
$ws.Range('J11').Value = '1005.3 hPa'
$ws.Range('E12').Value = '2026-02-07 04:18:02'
$ws.Range('H12').NumberFormat = "@"
$ws.Range('H12').Value = '70%'
$ws.Range('N12').Value = '8.1 °C 3:58 TU'
$ws.Range('O12').Value = '10.1 °C'
$ws.Range('E13').Value = '2026-02-07 04:18:04'
$ws.Range('E14').Value = '2026-02-07 04:18:06'
$ws.Range('H14').NumberFormat = "@"
$ws.Range('H14').Value = '80%'
$ws.Range('E15').Value = '2026-02-07 04:18:09'
$ws.Range('H15').NumberFormat = "@"
$ws.Range('H15').Value = '81%'
$ws.Range('J15').Value = '1001.2 hPa'
$ws.Range('N15').Value = '4.3 °C 3:50 TU'
$ws.Range('O15').Value = '7.2 °C'
$ws.Range('E16').Value = '2026-02-07 04:18:11'
$ws.Range('H16').NumberFormat = "@"
$ws.Range('H16').Value = '89%'
$ws.Range('N16').Value = '2.1 °C 3:34 TU'
$ws.Range('O16').Value = '3.3 °C'
$ws.Range('E17').Value = '2026-02-07 04:18:14'
$ws.Range('E18').Value = '2026-02-07 04:18:16'
$ws.Range('N18').Value = '-8.4 °C 3:53 TU'
$ws.Range('O18').Value = '-6.9 °C'
$ws.Range('E19').Value = '2026-02-07 04:18:18'
$ws.Range('J19').Value = '1005.6 hPa'
$ws.Range('N19').Value = '3.4 °C 3:52 TU'
$ws.Range('O19').Value = '4.8 °C'
$ws.Range('E20').Value = '2026-02-07 04:18:21'
$ws.Range('H20').NumberFormat = "@"
$ws.Range('H20').Value = '86%'
$ws.Range('N20').Value = '-5.1 °C 3:59 TU'
$ws.Range('O20').Value = '-4.4 °C'
$ws.Range('E21').Value = '2026-02-07 04:18:23'
$ws.Range('H21').NumberFormat = "@"
$ws.Range('H21').Value = '70%'
$ws.Range('J21').Value = '1001.2 hPa'
$ws.Range('N21').Value = '4.5 °C 3:59 TU'
$ws.Range('O21').Value = '7.9 °C'
$ws.Range('E22').Value = '2026-02-07 04:18:26'
$ws.Range('L22').Value = '11.2 km/h - 329º 3:35 TU'
$ws.Range('N22').Value = '4.2 °C 3:34 TU'
$ws.Range('O22').Value = '5.6 °C'
$ws.Range('E23').Value = '2026-02-07 04:18:28'
$ws.Range('H23').NumberFormat = "@"
$ws.Range('H23').Value = '99%'
$ws.Range('J23').Value = '1001.1 hPa'
$ws.Range('E24').Value = '2026-02-07 04:18:30'
$ws.Range('H24').NumberFormat = "@"
$ws.Range('H24').Value = '81%'
$ws.Range('N24').Value = '9.9 °C 3:36 TU'
$ws.Range('O24').Value = '10.4 °C'
$ws.Range('E25').Value = '2026-02-07 04:18:33'
$ws.Range('J25').Value = '1004.9 hPa'
$ws.Range('E26').Value = '2026-02-07 04:18:35'
$ws.Range('H26').NumberFormat = "@"
$ws.Range('H26').Value = '74%'
$ws.Range('N26').Value = '-2.5 °C 3:49 TU'
$ws.Range('E27').Value = '2026-02-07 04:18:38'
$ws.Range('J27').Value = '1000.9 hPa'
$ws.Range('L27').Value = '13.3 km/h - 290º 3:48 TU'
$ws.Range('E28').Value = '2026-02-07 04:18:40'
$ws.Range('H28').NumberFormat = "@"
$ws.Range('H28').Value = '86%'
$ws.Range('J28').Value = '1003.4 hPa'
$ws.Range('N28').Value = '2.0 °C 3:56 TU'
$ws.Range('O28').Value = '3.7 °C'
$ws.Range('E29').Value = '2026-02-07 04:18:42'
$ws.Range('H29').NumberFormat = "@"
$ws.Range('H29').Value = '57%'
$ws.Range('N29').Value = '9.4 °C 3:59 TU'
$ws.Range('O29').Value = '11.4 °C'
$ws.Range('E30').Value = '2026-02-07 04:18:44'
$ws.Range('E31').Value = '2026-02-07 04:18:47'
$ws.Range('J31').Value = '1005.5 hPa'
$ws.Range('O31').Value = '3.7 °C'
$ws.Range('E32').Value = '2026-02-07 04:18:49'
$ws.Range('J32').Value = '1003.9 hPa'
$ws.Range('E33').Value = '2026-02-07 04:18:52'
$ws.Range('N33').Value = '5.5 °C 3:56 TU'
$ws.Range('O33').Value = '7.4 °C'
$ws.Range('E34').Value = '2026-02-07 04:18:54'
$ws.Range('H34').NumberFormat = "@"
$ws.Range('H34').Value = '76%'
$ws.Range('N34').Value = '5.1 °C 3:59 TU'
$ws.Range('O34').Value = '6.7 °C'
$ws.Range('E35').Value = '2026-02-07 04:18:56'
$ws.Range('H35').NumberFormat = "@"
$ws.Range('H35').Value = '93%'
$ws.Range('N35').Value = '-7.7 °C 3:59 TU'
$ws.Range('O35').Value = '-4.7 °C'
$ws.Range('E36').Value = '2026-02-07 04:18:59'
$ws.Range('J36').Value = '1006.1 hPa'
$ws.Range('N36').Value = '4.1 °C 3:59 TU'
$ws.Range('O36').Value = '4.7 °C'
